# A new September SMS log entry ("sip") was recorded ahead of the most
# recent one, so row 48 on the "2024" sheet gets a new row inserted above
# it (pushing the existing September_Details/September_Date rows - and
# everything below them, including the trailing "Broadband" group label
# that lived at the bottom of the sheet - down by one row), and the new
# row is populated with the latest message/time pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("48:48").Insert()

$ws.Range("R48").Value = "sip"
$ws.Range("S48").Value = "2024-09-24 08:12:52"
